$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the bordered/bold ticker-label style (already used by A2:A7) down
# through the new rows (A8:A26) before writing values, so every label in the
# expanded table keeps the same look as the original six rows.
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A8:A26").PasteSpecial(-4122)

# (row, ticker, mean return, risk) - the table now points at the graph + the
# full stocks/crypto universe, sorted alphabetically by ticker.
$data = @(
    @(2, "AAPL", 0.000409195432119712, 0.01577304597325203),
    @(3, "ADA-USD", 0.002109322907709398, 0.05385911692947309),
    @(4, "AMZN", -0.00006331447102676818, 0.0197547657637857),
    @(5, "BABA", -0.0006555208053647794, 0.02917146136629337),
    @(6, "BNB-USD", 0.003300274737707548, 0.05275003261595376),
    @(7, "BTC-USD", 0.001338890310131067, 0.03414030084679871),
    @(8, "DAI-USD", -0.00002317790536334618, 0.002147598150660589),
    @(9, "DOGE-USD", 0.007118002616628462, 0.1316804669561277),
    @(10, "DOT-USD", 0.001307446046721384, 0.05807611652739767),
    @(11, "ETH-USD", 0.002288183890886724, 0.0450493043897679),
    @(12, "GOOG", 0.0005561434150984816, 0.01657202065415029),
    @(13, "GOOGL", 0.0005525675125263694, 0.01659933490063376),
    @(14, "JNJ", 0.00008346392399517444, 0.008663813227800198),
    @(15, "JPM", 0.0004363184681520799, 0.01386509989112363),
    @(16, "MATIC-USD", 0.005231289252545505, 0.07498958764790675),
    @(17, "MSFT", 0.0004038578094758204, 0.01497815621918233),
    @(18, "SOL-USD", 0.004044725884465609, 0.07182484297952729),
    @(19, "TON-USD", 0.0006277791916024, 0.06824254427809426),
    @(20, "TRX-USD", 0.002231566251630027, 0.04737623518475711),
    @(21, "USDC-USD", -0.000001439140654660258, 0.001241952029302097),
    @(22, "USDT-USD", -0.000002050810642745787, 0.0007281035275920827),
    @(23, "V", 0.000151766101611962, 0.01318653829784571),
    @(24, "VOD", -0.0002438401746817289, 0.014529915233304),
    @(25, "WBTC-USD", 0.001281691832933432, 0.03417820766614947),
    @(26, "XRP-USD", 0.002536448557910965, 0.06508533059648003)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

Write-Output "Updated A2:C26 with the graph/stocks return+risk table"
